$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-12-07 Sunday" "2025-12-08 Monday"

Replace-Text "445×7=" "153×3="
Replace-Text "453×2=" "701×5="
Replace-Text "807×3=" "499×4="
Replace-Text "576×7=" "120×9="
Replace-Text "736×2=" "534×9="

Replace-Text "748×7=" "493×2="
Replace-Text "735×6=" "271×5="
Replace-Text "990×5=" "942×3="
Replace-Text "744×2=" "659×4="
Replace-Text "882×5=" "424×8="

Replace-Text "438×2=" "502×3="
Replace-Text "977×8=" "408×8="
Replace-Text "435×8=" "510×8="
Replace-Text "489×2=" "785×5="
Replace-Text "911×9=" "127×8="

Replace-Text "565×6=" "310×9="
Replace-Text "481×4=" "608×7="
Replace-Text "198×8=" "188×7="
Replace-Text "851×9=" "600×4="
Replace-Text "168×5=" "930×6="

Replace-Text "784×3=" "878×9="
Replace-Text "707×9=" "524×6="
Replace-Text "954×7=" "584×2="
Replace-Text "677×4=" "648×3="
Replace-Text "558×3=" "261×7="
